# Fruta / hortaliza, semanal
# Re-shuffle the weekly rows (2..12) of the Papaya sheet: each row keeps its
# fixed/common columns, but the week-specific data (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Precio $/Kg) is reassigned to a new date row
# per the updated weekly report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "variable" columns (D, L, M, N, O, P, S) for rows 2..12
# before any writes, keyed by the original row number.
$orig = @{}
for ($r = 2; $r -le 12; $r++) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Destination row -> source row (which original row's data now lands here).
$mapping = @{
    2  = 10
    3  = 7
    4  = 11
    5  = 12
    6  = 3
    7  = 4
    8  = 6
    9  = 8
    10 = 2
    11 = 5
    12 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $orig[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2 = $data.D
    $ws.Cells.Item($destRow, 12).Value2 = $data.L
    $ws.Cells.Item($destRow, 13).Value2 = $data.M
    $ws.Cells.Item($destRow, 14).Value2 = $data.N
    $ws.Cells.Item($destRow, 15).Value2 = $data.O
    $ws.Cells.Item($destRow, 16).Value2 = $data.P
    $ws.Cells.Item($destRow, 19).Value2 = $data.S
}
